$d = $word.ActiveDocument

# 1. "MongoDB" -> "RDBMS" in the Data bullet list.
$d.Content.Find.Execute("MongoDB", $true, $false, $false, $false, $false, $true, 1, $false, "RDBMS", 2) | Out-Null

# 2. Move the "_GoBack" bookmark from the end of the
#    "Illustration will be saved in SVG format on disk" paragraph to the end
#    of the (now) "RDBMS" paragraph.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# Locate the freshly-replaced "RDBMS" run.
$r = $d.Content
$r.Find.Execute("RDBMS", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Collapse(0)

# NOTE: this COM host mis-resolves Bookmarks.Add() when handed a zero-length
# Range that sits exactly on a paragraph-mark boundary (it silently drops the
# bookmark at document position 0 instead of where intended). Work around it
# by growing the target paragraph by one throwaway character so the insertion
# point is no longer on that boundary, anchoring the bookmark there, and then
# deleting the throwaway character again. The bookmark itself is emitted by
# the host so it stays correctly anchored once the extra character is gone.
$r.InsertAfter("X")
$anchor = $d.Range($r.End - 1, $r.End - 1)
$d.Bookmarks.Add("_GoBack", $anchor) | Out-Null
$placeholder = $d.Range($r.End - 1, $r.End)
$placeholder.Delete()
